# Auto-generated Excel COM-interop script to apply scraped price updates
# across the 8 job-leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 4166789.2
$ws.Range("I33").Value = 4545581.5
$ws.Range("K33").Value = 4545581.5
$ws.Range("M33").Value = -4545352.5
# Row 100
$ws.Range("H100").Value = 2072.8096
$ws.Range("I100").Value = 2243.3333
$ws.Range("K100").Value = 2243.3333
$ws.Range("M100").Value = -1702.3333
# Row 115
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = $null
$ws.Range("N115").Value = $null
# Row 116
$ws.Range("H116").Value = 12036.6
$ws.Range("I116").Value = 11706.667
$ws.Range("K116").Value = 11706.667
$ws.Range("M116").Value = -8264.666999999999
# Row 137
$ws.Range("H137").Value = 1285.25
$ws.Range("I137").Value = 1041.5294
$ws.Range("J137").Value = 2666.3333
$ws.Range("K137").Value = 3124.5882
$ws.Range("L137").Value = 7998.999899999999
$ws.Range("M137").Value = -574.5881999999997
$ws.Range("N137").Value = -13098.9999
# Row 141
$ws.Range("H141").Value = 300450
$ws.Range("I141").Value = 900
$ws.Range("J141").Value = 600000
$ws.Range("K141").Value = 2700
$ws.Range("L141").Value = 1800000
$ws.Range("M141").Value = 2480
$ws.Range("N141").Value = -1810360

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7578.8823
$ws.Range("I32").Value = 7957
$ws.Range("K32").Value = 7957
$ws.Range("M32").Value = -7670
# Row 97
$ws.Range("H97").Value = 1399.0869
$ws.Range("I97").Value = 1431.8096
$ws.Range("J97").Value = 1055.5
$ws.Range("K97").Value = 1431.8096
$ws.Range("L97").Value = 1055.5
$ws.Range("M97").Value = -935.8096
$ws.Range("N97").Value = -2047.5
# Row 122
$ws.Range("H122").Value = 2124.1667
$ws.Range("I122").Value = 2049.1177
$ws.Range("K122").Value = 6147.353099999999
$ws.Range("M122").Value = -3697.353099999999
# Row 132
$ws.Range("H132").Value = 3476.6956
$ws.Range("I132").Value = 3331.6667
$ws.Range("K132").Value = 9995.000100000001
$ws.Range("M132").Value = -7465.000100000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 5324.4287
$ws.Range("I20").Value = 5259.913
$ws.Range("J20").Value = 5448.0835
$ws.Range("K20").Value = 5259.913
$ws.Range("L20").Value = 5448.0835
$ws.Range("M20").Value = -5012.913
$ws.Range("N20").Value = -5942.0835
# Row 86
$ws.Range("H86").Value = 12120.477
$ws.Range("I86").Value = 11148.412
$ws.Range("J86").Value = 16251.75
$ws.Range("K86").Value = 11148.412
$ws.Range("L86").Value = 16251.75
$ws.Range("M86").Value = -10025.412
$ws.Range("N86").Value = -18497.75
# Row 89
$ws.Range("H89").Value = 12120.477
$ws.Range("I89").Value = 11148.412
$ws.Range("J89").Value = 16251.75
$ws.Range("K89").Value = 55742.06
$ws.Range("L89").Value = 81258.75
$ws.Range("M89").Value = -50126.06
$ws.Range("N89").Value = -92490.75
# Row 105
$ws.Range("H105").Value = 2505.0312
$ws.Range("I105").Value = 2467.7917
$ws.Range("J105").Value = 2616.75
$ws.Range("K105").Value = 2467.7917
$ws.Range("L105").Value = 2616.75
$ws.Range("M105").Value = -720.7917000000002
$ws.Range("N105").Value = -6110.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 77120.91
$ws.Range("J31").Value = 19966.666
$ws.Range("L31").Value = 19966.666
$ws.Range("N31").Value = -20556.666
# Row 34
$ws.Range("H34").Value = 77120.91
$ws.Range("J34").Value = 19966.666
$ws.Range("L34").Value = 19966.666
$ws.Range("N34").Value = -20370.666
# Row 62
$ws.Range("H62").Value = 5467
$ws.Range("I62").Value = 5509.4546
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 5509.4546
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -4885.4546
$ws.Range("N62").Value = -6248
# Row 65
$ws.Range("H65").Value = 5467
$ws.Range("I65").Value = 5509.4546
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 27547.273
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -24427.273
$ws.Range("N65").Value = -31240
# Row 86
$ws.Range("H86").Value = 5092.143
$ws.Range("I86").Value = 5092.143
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 5092.143
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -3969.143
$ws.Range("N86").Value = $null
# Row 88
$ws.Range("H88").Value = 26895
$ws.Range("J88").Value = 28474.166
$ws.Range("L88").Value = 28474.166
$ws.Range("N88").Value = -29286.166
# Row 89
$ws.Range("H89").Value = 5092.143
$ws.Range("I89").Value = 5092.143
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 25460.715
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -19844.715
$ws.Range("N89").Value = $null
# Row 91
$ws.Range("H91").Value = 26895
$ws.Range("J91").Value = 28474.166
$ws.Range("L91").Value = 28474.166
$ws.Range("N91").Value = -31282.166
# Row 105
$ws.Range("H105").Value = 1686.75
$ws.Range("J105").Value = 2366.6667
$ws.Range("L105").Value = 2366.6667
$ws.Range("N105").Value = -5860.6667
# Row 134
$ws.Range("H134").Value = 11461.5
$ws.Range("I134").Value = 5296.552
$ws.Range("J134").Value = 37002
$ws.Range("K134").Value = 15889.656
$ws.Range("L134").Value = 111006
$ws.Range("M134").Value = -13354.656
$ws.Range("N134").Value = -116076

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 71377.875
$ws.Range("I131").Value = 182400.9
$ws.Range("J131").Value = 13222.952
$ws.Range("K131").Value = 547202.7
$ws.Range("L131").Value = 39668.856
$ws.Range("M131").Value = -542162.7
$ws.Range("N131").Value = -49748.856

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 10671.272
$ws.Range("J70").Value = 7333.3335
$ws.Range("L70").Value = 7333.3335
$ws.Range("N70").Value = -7873.3335
# Row 73
$ws.Range("H73").Value = 10671.272
$ws.Range("J73").Value = 7333.3335
$ws.Range("L73").Value = 7333.3335
$ws.Range("N73").Value = -9205.333500000001
# Row 80
$ws.Range("H80").Value = 2833.3333
$ws.Range("I80").Value = 2750
$ws.Range("K80").Value = 2750
$ws.Range("M80").Value = -1752
# Row 83
$ws.Range("H83").Value = 2833.3333
$ws.Range("I83").Value = 2750
$ws.Range("K83").Value = 13750
$ws.Range("M83").Value = -8758
# Row 126
$ws.Range("H126").Value = 32182.637
$ws.Range("I126").Value = 48272.715
$ws.Range("J126").Value = 4025
$ws.Range("K126").Value = 144818.145
$ws.Range("L126").Value = 12075
$ws.Range("M126").Value = -142348.145
$ws.Range("N126").Value = -17015

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 23
$ws.Range("H23").Value = 3250000
$ws.Range("I23").Value = 1000000
$ws.Range("J23").Value = 5500000
$ws.Range("K23").Value = 1000000
$ws.Range("L23").Value = 5500000
$ws.Range("M23").Value = -999770
$ws.Range("N23").Value = -5500460
# Row 46
$ws.Range("H46").Value = 1691.3334
$ws.Range("I46").Value = 1850
$ws.Range("J46").Value = 1374
$ws.Range("K46").Value = 1850
$ws.Range("L46").Value = 1374
$ws.Range("M46").Value = -1662
$ws.Range("N46").Value = -1750
# Row 136
$ws.Range("H136").Value = 6592.3
$ws.Range("I136").Value = 5991
$ws.Range("K136").Value = 17973
$ws.Range("M136").Value = -15423

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 1845.0714
$ws.Range("I136").Value = 1648.3334
$ws.Range("K136").Value = 4945.0002
$ws.Range("M136").Value = -2395.0002

